$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.279.58"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "1.898.65"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'323.91"
$ws.Range("E5").Value = "  -1.82%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.4710"
$ws.Range("E7").Value = "  +2.88%  "

$ws.Range("D8").Value = "'0.4027"
$ws.Range("E8").Value = "  -1.85%  "

$ws.Range("D9").Value = "'47.66"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "'0.08004"
$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("D11").Value = "'0.9952"
$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("D12").Value = "'22.56"
$ws.Range("E12").Value = "  +3.88%  "

$ws.Range("D13").Value = "1.907.92"
$ws.Range("E13").Value = "  +0.29%  "

$ws.Range("D14").Value = "'5.862"
$ws.Range("E14").Value = "  -0.95%  "

$ws.Range("D15").Value = "'7.048"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").Value = "'89.06"
$ws.Range("E16").Value = "  +0.45%  "

$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").Value = "'0.06619"
$ws.Range("E18").Value = "  +1.03%  "

$ws.Range("D19").Value = "'0.00001027"
$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").Value = "'17.50"
$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("D22").Value = "29.242.78"
$ws.Range("E22").Value = "  +0.69%  "

$ws.Range("D23").Value = "'5.502"
$ws.Range("E23").Value = "  +1.33%  "

$ws.Range("E24").Value = "  +1.18%  "

$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("D26").Value = "2.113.03"
$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("D27").Value = "'154.31"
$ws.Range("E27").Value = "  -0.95%  "

$ws.Range("D28").Value = "'19.68"
$ws.Range("E28").Value = "  +0.56%  "

$ws.Range("D29").Value = "'6.087"
$ws.Range("E29").Value = "  +10.74%  "

$ws.Range("D30").Value = "'2.089"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").Value = "'117.37"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").Value = "'1.061"
$ws.Range("E32").Value = "  +2.44%  "

$ws.Range("D33").Value = "'0.09454"
$ws.Range("E33").Value = "  +1.55%  "

$ws.Range("D34").Value = "'1.400"
$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("D35").Value = "'3.543"
$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").Value = "'5.349"
$ws.Range("E36").Value = "  +1.06%  "

$ws.Range("D37").Value = "'0.06085"
$ws.Range("E37").Value = "  +0.57%  "

$ws.Range("D38").Value = "'0.02245"
$ws.Range("E38").Value = "  +0.82%  "

$ws.Range("D39").Value = "'1.172"
$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").Value = "'8.078"
$ws.Range("E40").Value = "  -3.35%  "

$ws.Range("D41").Value = "'0.5809"
$ws.Range("E41").Value = "  +0.45%  "

$ws.Range("D42").Value = "'0.1829"
$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("D43").Value = "'2.478"
$ws.Range("E43").Value = "  +8.72%  "

$ws.Range("D44").Value = "'10.06"
$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("E45").Value = "  +0.71%  "

$ws.Range("D46").Value = "'0.07707"
$ws.Range("E46").Value = "  +2.62%  "

$ws.Range("D47").Value = "'12.08"
$ws.Range("E47").Value = "  +1.13%  "

$ws.Range("D48").Value = "'0.5481"
$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("D49").Value = "'1.901"
$ws.Range("E49").Value = "  -0.06%  "

$ws.Range("D50").Value = "'113.45"

$ws.Range("D51").Value = "'43.97"
$ws.Range("E51").Value = "  -0.90%  "
